$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "202.173.124.126"
$ws.Range("B12").Value = 28.3621533
$ws.Range("C12").Value = 77.2828576
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = "Mozilla/5.0 (Linux; Android 10; K) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Mobile Safari/537.36"
$ws.Range("F12").Value = "Linux armv81"
$ws.Range("G12").Value = "2025-06-25T16:50:01.397Z"
